$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell -> new value pairs scraped from the crypto-ticker refresh.
# Values that look like plain numbers (e.g. "2.33") must be forced to
# Text so Excel does not silently coerce them (stripping the trailing
# zero in "0.110", turning "81.44" into a float, etc.) - matching the
# source data which stores these as inline strings.
$updates = @(
    @{Cell="D2"; Value="63.949.48"}
    @{Cell="E2"; Value="  +0.27%  "}
    @{Cell="D3"; Value="3.134.35"}
    @{Cell="E3"; Value="  +0.63%  "}
    @{Cell="E4"; Value="  +0.02%  "}
    @{Cell="D5"; Value="589.46"}
    @{Cell="E5"; Value="  +0.69%  "}
    @{Cell="D6"; Value="145.05"}
    @{Cell="E7"; Value="  +0.00%  "}
    @{Cell="D8"; Value="3.128.06"}
    @{Cell="E8"; Value="  +0.64%  "}
    @{Cell="E9"; Value="  -0.14%  "}
    @{Cell="E10"; Value="  -0.32%  "}
    @{Cell="D11"; Value="5.92"}
    @{Cell="E11"; Value="  +2.77%  "}
    @{Cell="E12"; Value="  -1.89%  "}
    @{Cell="E13"; Value="  -2.20%  "}
    @{Cell="D14"; Value="37.35"}
    @{Cell="E14"; Value="  +0.52%  "}
    @{Cell="D15"; Value="3.653.54"}
    @{Cell="E15"; Value="  +0.68%  "}
    @{Cell="E16"; Value="  -1.24%  "}
    @{Cell="D17"; Value="7.33"}
    @{Cell="E17"; Value="  +2.83%  "}
    @{Cell="D18"; Value="63.772.61"}
    @{Cell="E18"; Value="  +0.12%  "}
    @{Cell="D19"; Value="3.131.19"}
    @{Cell="E19"; Value="  +0.47%  "}
    @{Cell="D20"; Value="466.61"}
    @{Cell="E20"; Value="  +0.40%  "}
    @{Cell="D21"; Value="14.34"}
    @{Cell="E21"; Value="  +0.26%  "}
    @{Cell="E22"; Value="  +0.17%  "}
    @{Cell="D23"; Value="7.53"}
    @{Cell="E23"; Value="  +0.23%  "}
    @{Cell="B24"; Value="Fetch.AI"}
    @{Cell="C24"; Value="https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"}
    @{Cell="D24"; Value="2.33"}
    @{Cell="E24"; Value="  +7.98%  "}
    @{Cell="D25"; Value="12.96"}
    @{Cell="E25"; Value="  -1.18%  "}
    @{Cell="B26"; Value="Litecoin"}
    @{Cell="C26"; Value="https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"}
    @{Cell="D26"; Value="81.44"}
    @{Cell="E26"; Value="  -0.53%  "}
    @{Cell="E27"; Value="  +0.05%  "}
    @{Cell="D28"; Value="9.92"}
    @{Cell="E28"; Value="  +10.91%  "}
    @{Cell="D29"; Value="7.46"}
    @{Cell="E29"; Value="  +8.96%  "}
    @{Cell="D30"; Value="2.71"}
    @{Cell="E30"; Value="  +0.47%  "}
    @{Cell="D31"; Value="2.24"}
    @{Cell="E31"; Value="  +0.37%  "}
    @{Cell="E32"; Value="  +0.07%  "}
    @{Cell="D33"; Value="27.64"}
    @{Cell="E33"; Value="  +2.63%  "}
    @{Cell="D34"; Value="0.110"}
    @{Cell="E34"; Value="  +0.64%  "}
    @{Cell="E35"; Value="  -2.69%  "}
    @{Cell="E36"; Value="  +0.68%  "}
    @{Cell="E37"; Value="  +1.11%  "}
    @{Cell="E38"; Value="  -2.30%  "}
    @{Cell="E39"; Value="  -6.03%  "}
    @{Cell="D40"; Value="51.32"}
    @{Cell="E40"; Value="  +0.76%  "}
    @{Cell="D41"; Value="9.35"}
    @{Cell="E41"; Value="  +7.86%  "}
    @{Cell="D42"; Value="452.51"}
    @{Cell="E42"; Value="  +1.20%  "}
    @{Cell="E43"; Value="  +5.29%  "}
    @{Cell="E44"; Value="  +0.31%  "}
    @{Cell="D45"; Value="2.912.79"}
    @{Cell="E45"; Value="  +1.24%  "}
    @{Cell="D46"; Value="40.24"}
    @{Cell="E46"; Value="  +12.47%  "}
    @{Cell="D47"; Value="0.107"}
    @{Cell="E47"; Value="  -3.03%  "}
    @{Cell="D48"; Value="133.16"}
    @{Cell="E48"; Value="  +7.86%  "}
    @{Cell="E49"; Value="  -0.02%  "}
    @{Cell="B50"; Value="ThetaToken"}
    @{Cell="C50"; Value="https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"}
    @{Cell="D50"; Value="2.24"}
    @{Cell="E50"; Value="  +2.53%  "}
    @{Cell="B51"; Value="Stellar"}
    @{Cell="C51"; Value="https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"}
    @{Cell="D51"; Value="0.111"}
    @{Cell="E51"; Value="  -0.58%  "}
)

foreach ($u in $updates) {
    $cell = $ws.Range($u.Cell)
    $v = $u.Value
    $looksNumeric = $v -match "^\s*[+-]?(\d+\.?\d*|\.\d+)\s*$"
    if ($looksNumeric) {
        # Force text storage so Excel keeps the exact literal digits
        $cell.NumberFormat = "@"
        $cell.Value = $v
        $cell.ClearFormats()
    } else {
        $cell.Value = $v
    }
}
